# Re-order the comma-separated "Recorded By" names/emails in column G so
# that "System" (exact, case-sensitive match) comes first if present,
# otherwise "admin@admin.com" comes first if present. The relative order
# of the remaining entries is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ContainsExact($arr, $val) {
    foreach ($x in $arr) {
        if ($x.Equals($val)) { return $true }
    }
    return $false
}

function ReorderRecordedBy($text) {
    $parts = $text.Split(",")
    $items = @()
    foreach ($p in $parts) {
        $items += $p.Trim()
    }

    if ($items.Count -le 1) {
        return $text
    }

    $priorityValue = $null
    if (ContainsExact $items "System") {
        $priorityValue = "System"
    } elseif (ContainsExact $items "admin@admin.com") {
        $priorityValue = "admin@admin.com"
    } else {
        return $text
    }

    $rest = @()
    $removed = $false
    foreach ($it in $items) {
        if ((-not $removed) -and $it.Equals($priorityValue)) {
            $removed = $true
            continue
        }
        $rest += $it
    }

    $newItems = @($priorityValue) + $rest
    return [string]::Join(", ", $newItems)
}

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $newText = ReorderRecordedBy $text

    if (-not $newText.Equals($text)) {
        $cell.Value = $newText
    }
}
